$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: update problem/solution description, remove Developer (Viet)
$ws.Range("B2").Value = "Lỗi tạo ra con trỏ 'Pointer.cur'"
$ws.Range("C2").Value = "Không xuất hiện lỗi"
$ws.Range("E2").Clear()

# Row 3: clear out the old issue entirely (keep B3's style, but no content)
$ws.Range("B3").ClearContents()
$ws.Range("C3").Clear()
$ws.Range("D3").Clear()

# Row 4: clear out the old issue entirely (keep B4's style, but no content)
$ws.Range("B4").ClearContents()
$ws.Range("C4").Clear()
$ws.Range("D4").Clear()

# Row 5: clear out the old issue entirely
$ws.Range("B5").Clear()
$ws.Range("C5").Clear()
$ws.Range("D5").Clear()

# Update selection to C2
$ws.Range("C2").Select()
